$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New calculation_time values (column F) for rows 2-115,
# reflecting the metrics recalculation timestamp update.
$newTimes = @(
    45774.51714972244,
    45774.51715049539,
    45774.5171510482,
    45774.51715155634,
    45774.51715210335,
    45774.51715266592,
    45774.51715319938,
    45774.51715377426,
    45774.51715433363,
    45774.51715486211,
    45774.51715539256,
    45774.51715595108,
    45774.51715651037,
    45774.51715703134,
    45774.51715754552,
    45774.51715805998,
    45774.51715857739,
    45774.51715910548,
    45774.51715972563,
    45774.5171603612,
    45774.51716100444,
    45774.51716164342,
    45774.51716231722,
    45774.51716316248,
    45774.51716414197,
    45774.51716510336,
    45774.5171659611,
    45774.51716695887,
    45774.51716788067,
    45774.51716863063,
    45774.51716917584,
    45774.51716970009,
    45774.51717024361,
    45774.51717077286,
    45774.51717126879,
    45774.51717177722,
    45774.5171722889,
    45774.51717278786,
    45774.51717330979,
    45774.51717382465,
    45774.5171743219,
    45774.51717483098,
    45774.51717533365,
    45774.5171758264,
    45774.51717632551,
    45774.51717681944,
    45774.51717731839,
    45774.51717788205,
    45774.51717849063,
    45774.51717910669,
    45774.5171796872,
    45774.51718033342,
    45774.51718098157,
    45774.51718151278,
    45774.51718201656,
    45774.51718251223,
    45774.51718302231,
    45774.51718302244,
    45774.51718373374,
    45774.51718454775,
    45774.51718513152,
    45774.51718513161,
    45774.51718513163,
    45774.51718513164,
    45774.51718513166,
    45774.51718513167,
    45774.51718513168,
    45774.51718513175,
    45774.51718513176,
    45774.51718513177,
    45774.51718513178,
    45774.51718513179,
    45774.51718566574,
    45774.51718566582,
    45774.51718566583,
    45774.51718566584,
    45774.51718566586,
    45774.51718566587,
    45774.51718566588,
    45774.51718566588,
    45774.51718566589,
    45774.5171856659,
    45774.51718566591,
    45774.51718566593,
    45774.51718689167,
    45774.51718746465,
    45774.5171884009,
    45774.51718842874,
    45774.51718845382,
    45774.51718847759,
    45774.51718850067,
    45774.51718852334,
    45774.51718854681,
    45774.51718856983,
    45774.51718859068,
    45774.51718861112,
    45774.51718863122,
    45774.51718865108,
    45774.51718867097,
    45774.51718869057,
    45774.5171887101,
    45774.51718873045,
    45774.51718874928,
    45774.51718876838,
    45774.51718878747,
    45774.51718880704,
    45774.5171888262,
    45774.51718884545,
    45774.51718886448,
    45774.51718888384,
    45774.5171889031,
    45774.51718892213,
    45774.51718894127,
    45774.51718896026
)

$startRow = 2
for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 6).Value = $newTimes[$i]
}

Write-Host "Updated calculation_time for $($newTimes.Length) rows"
